$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking count (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update total correct count (B12): 57 -> 95
$ws.Range("B12").Value = 95

# Update correct/total marks label (E12): "54/84" -> "95/140"
$ws.Range("E12").Value = "95/140"
